# Update the EPEX Spot prices workbook with the newest day of data.
#
# Sheet "Prix Spot": add a new date column (CD) "03-sep" with 24 hourly values.
# Sheet "Gaz":       add a new row (79) for date 2025-09-01 / 30.8.
# Sheet "CO2":       add a new row (79) for date 2025-09-01 / 73.31.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Prix Spot" - new column CD ("03-sep")
# ---------------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Give CD1 the same look (bold / centered / bordered) as the rest of row 1
# by copying the formatting of the previous header cell (CC1) first.
$wsSpot.Range("CC1").Copy()
$wsSpot.Range("CD1").PasteSpecial(-4122)   # xlPasteFormats
$wsSpot.Range("CD1").Value = "03-sep"

$spotValues = @(
    17.44,
    15.65,
    10.86,
    8.460000000000001,
    4.11,
    8.94,
    10.14,
    21.34,
    22.6,
    16.37,
    0,
    -0.01,
    -0.01,
    -0.02,
    -0.02,
    -0.01,
    -0.01,
    0,
    12.85,
    56.79,
    58.21,
    53.41,
    70.23,
    56.33
)

for ($i = 0; $i -lt $spotValues.Length; $i++) {
    $row = $i + 2
    $wsSpot.Cells.Item($row, 82).Value = $spotValues[$i]
}

# ---------------------------------------------------------------------------
# Sheet 2: "Gaz" - new row 79 (2025-09-01 / 30.8)
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the date to be stored as literal text (matching every other date
# cell in column A) instead of letting Excel auto-convert it to a serial
# date number.
$wsGaz.Range("A79").NumberFormat = "@"
$wsGaz.Range("A79").Value = "2025-09-01"
# Reset the cell's look back to the plain (unformatted) style used by the
# rest of column A.
$wsGaz.Range("A78").Copy()
$wsGaz.Range("A79").PasteSpecial(-4122)    # xlPasteFormats

$wsGaz.Range("B79").Value = 30.8

# ---------------------------------------------------------------------------
# Sheet 3: "CO2" - new row 79 (2025-09-01 / 73.31)
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A79").NumberFormat = "@"
$wsCo2.Range("A79").Value = "2025-09-01"
$wsCo2.Range("A78").Copy()
$wsCo2.Range("A79").PasteSpecial(-4122)    # xlPasteFormats

$wsCo2.Range("B79").Value = 73.31
